# Update column G ("K") values for rows 2-30 of Sheet1 in save_data workbook.
# This regenerates the K column (previously "Strike#") with recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 2
    6  = 3
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 3
    20 = 2
    21 = 1
    22 = 1
    23 = 2
    24 = 3
    25 = 1
    26 = 1
    27 = 2
    28 = 2
    29 = 1
    30 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
